# Applies the "/diff" column ('correct' marker in column C) to the three
# défi5 vocabulary sheets, plus a handful of French/Dutch text typo fixes
# that were corrected in the same commit.

$wb = $excel.ActiveWorkbook

# --- défi5m1 (sheet index 13): 51 rows ---------------------------------
$ws1 = $wb.Worksheets.Item(13)

# Typo fixes in the source/target text columns.
$ws1.Cells.Item(3, 1).Value  = "un besoin financier"
$ws1.Cells.Item(6, 1).Value  = "un commerçant du coin, une commerçante du coin"
$ws1.Cells.Item(19, 1).Value = "dépendant, dépendante [de]"
$ws1.Cells.Item(19, 2).Value = "afhankelijk [van]"
$ws1.Cells.Item(29, 1).Value = "contenir"

for ($r = 1; $r -le 51; $r++) {
    $ws1.Cells.Item($r, 3).Value = "correct"
}

# --- défi5m2 (sheet index 14): 54 rows ---------------------------------
$ws2 = $wb.Worksheets.Item(14)

$ws2.Cells.Item(1, 1).Value  = "un aide-magasin, une aide-magasin"
$ws2.Cells.Item(47, 1).Value = "à temps plein"

for ($r = 1; $r -le 54; $r++) {
    $ws2.Cells.Item($r, 3).Value = "correct"
}

# --- défi5m3 (sheet index 15): 19 rows ---------------------------------
$ws3 = $wb.Worksheets.Item(15)

for ($r = 1; $r -le 19; $r++) {
    $ws3.Cells.Item($r, 3).Value = "correct"
}
